# Update "想去人数" (F column) and one "最低票价" (G column) figures
# across the "展览" and "全部类型" worksheets, per updated scrape data.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, column, newValue) updates.
$updates = @{
    "展览" = @(
        @{ Row = 2;  Col = "F"; Value = 6967 },
        @{ Row = 4;  Col = "F"; Value = 455 },
        @{ Row = 7;  Col = "F"; Value = 127 },
        @{ Row = 8;  Col = "F"; Value = 115 },
        @{ Row = 11; Col = "G"; Value = 45 },
        @{ Row = 13; Col = "F"; Value = 442 },
        @{ Row = 15; Col = "F"; Value = 1815 },
        @{ Row = 16; Col = "F"; Value = 40 },
        @{ Row = 17; Col = "F"; Value = 3579 },
        @{ Row = 20; Col = "F"; Value = 82 },
        @{ Row = 21; Col = "F"; Value = 15 },
        @{ Row = 23; Col = "F"; Value = 2198 },
        @{ Row = 24; Col = "F"; Value = 2 },
        @{ Row = 25; Col = "F"; Value = 231 },
        @{ Row = 27; Col = "F"; Value = 33 },
        @{ Row = 31; Col = "F"; Value = 150 },
        @{ Row = 32; Col = "F"; Value = 108 },
        @{ Row = 33; Col = "F"; Value = 53 }
    )
    "全部类型" = @(
        @{ Row = 2;  Col = "F"; Value = 6967 },
        @{ Row = 4;  Col = "F"; Value = 455 },
        @{ Row = 8;  Col = "F"; Value = 127 },
        @{ Row = 9;  Col = "F"; Value = 115 },
        @{ Row = 12; Col = "G"; Value = 45 },
        @{ Row = 14; Col = "F"; Value = 442 },
        @{ Row = 16; Col = "F"; Value = 1815 },
        @{ Row = 17; Col = "F"; Value = 40 },
        @{ Row = 18; Col = "F"; Value = 3579 },
        @{ Row = 21; Col = "F"; Value = 82 },
        @{ Row = 22; Col = "F"; Value = 15 },
        @{ Row = 24; Col = "F"; Value = 2198 },
        @{ Row = 25; Col = "F"; Value = 2 },
        @{ Row = 26; Col = "F"; Value = 231 },
        @{ Row = 28; Col = "F"; Value = 33 },
        @{ Row = 32; Col = "F"; Value = 150 },
        @{ Row = 33; Col = "F"; Value = 108 },
        @{ Row = 34; Col = "F"; Value = 53 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($change in $updates[$sheetName]) {
        $cellRef = "$($change.Col)$($change.Row)"
        $ws.Range($cellRef).Value = $change.Value
    }
}
